$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of existing death counts (column C) ---
$ws.Cells.Item(823, 3).Value = 8
$ws.Cells.Item(870, 3).Value = 15
$ws.Cells.Item(951, 3).Value = 28
$ws.Cells.Item(971, 3).Value = 26
$ws.Cells.Item(979, 3).Value = 14
$ws.Cells.Item(980, 3).Value = 32
$ws.Cells.Item(983, 3).Value = 32

# --- Rows 984-986 shift down one age group (values for date 44159) ---
$ws.Cells.Item(984, 2).Value = "50-59"
$ws.Cells.Item(984, 3).Value = 1

$ws.Cells.Item(985, 2).Value = "60-69"
$ws.Cells.Item(985, 3).Value = 6

$ws.Cells.Item(986, 2).Value = "70-79"
$ws.Cells.Item(986, 3).Value = 16

# --- Insert 10 new rows (987-996) for the new trailing data ---
$ws.Range("A987:A996").EntireRow.Insert()

$newRows = @(
    @(44159, "80+",   27),
    @(44160, "40-49",  1),
    @(44160, "50-59",  1),
    @(44160, "60-69",  4),
    @(44160, "70-79",  7),
    @(44160, "80+",   16),
    @(44161, "50-59",  2),
    @(44161, "60-69",  5),
    @(44161, "70-79",  4),
    @(44161, "80+",   15)
)

$r = 987
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
